# Generate Report for Handback
#
# Refreshes the handback-status report timestamps for the
# "36038c19-52ce-4a1b-8036-de19daaeacb8" source file's latest handoff/
# handback cycle across the Overview, zh-cn and de-de sheets. The second
# data row (fbebf677-...) already reflects an in-sync handback and is left
# untouched.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-23 00:45:57"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-23 00:45:52"
$zhcn.Range("K2").Value = "2016-08-23 00:46:14"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-23 00:45:57"
$dede.Range("K2").Value = "2016-08-23 00:46:21"
